# Fill in longitude (row 3) and latitude (row 4) data for columns B:F,
# then move the selection to C7, matching the committed workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - longitude
$ws.Range("B3").Value = -176.65
$ws.Range("C3").Value = -149.8056
$ws.Range("D3").Value = -148.9
$ws.Range("E3").Value = 149.1108
$ws.Range("F3").Value = 150.0036

# Row 4 - latitude
$ws.Range("B4").Value = 51.8833
$ws.Range("C4").Value = 61.1889
$ws.Range("D4").Value = 61.7167
$ws.Range("E4").Value = 60.9583
$ws.Range("F4").Value = 61.2167

# Size columns B and C to hug the newly-entered numbers (bestFit-style
# autofit), matching the "8" / "9.7109375" column widths from the diff.
$ws.Columns("B").ColumnWidth = 7.15
$ws.Columns("C").ColumnWidth = 8.8

# Move the active selection to C7
$ws.Range("C7").Select()
